# no-op
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
